# Insert one new weekly price-report row right before the existing row 380.
# This pushes the former rows 380-459 down to 381-460 (dimension grows from
# R459 to R460) and populates the newly opened row 380 with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 380..459 down by one to make room for the new record.
$ws.Rows.Item(380).Insert()

# Populate the newly inserted row 380 with the new weekly data point.
$ws.Range("A380").Value = 8
$ws.Range("B380").Value = "Terminal La Palmera de La Serena"
$ws.Range("C380").Value = "Coquimbo"
$ws.Range("D380").Value = 45209
$ws.Range("E380").Value = 4
$ws.Range("F380").Value = 100112031
$ws.Range("G380").Value = "Poroto verde"
$ws.Range("H380").Value = "Magnum"
$ws.Range("I380").Value = "Primera"
$ws.Range("J380").Value = 560
$ws.Range("K380").Value = 26000
$ws.Range("L380").Value = 27000
$ws.Range("M380").Value = 26500
$ws.Range("N380").Value = "$/malla 25 kilos"
$ws.Range("O380").Value = "Perú"
$ws.Range("P380").Value = 1060
$ws.Range("Q380").Value = 25
$ws.Range("R380").Value = "Hortaliza"
